$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows sourced from DGS's 2021/10/20 and 2021/10/22 reports.
# Column A holds the report date as text (displayed via the column's
# yyyy/mm/dd number format), so force a text format before assigning the
# value to avoid Excel auto-converting the string into a date serial,
# then restore the date display format used by the rest of the column.
$ws.Range("A96").NumberFormat = "@"
$ws.Range("A96").Value = "2021/10/20"
$ws.Range("A96").NumberFormat = "yyyy/mm/dd"

$ws.Range("B96").Value = 84.4
$ws.Range("C96").Value = 84.8
$ws.Range("D96").Value = 1.02
$ws.Range("E96").Value = 1.02

$ws.Range("A97").NumberFormat = "@"
$ws.Range("A97").Value = "2021/10/22"
$ws.Range("A97").NumberFormat = "yyyy/mm/dd"

$ws.Range("B97").Value = 86.1
$ws.Range("C97").Value = 86.5
$ws.Range("D97").Value = 1.02
$ws.Range("E97").Value = 1.02

$ws.Range("A98").Select()
